# Lesson 282 - Monday
# Reveal the missing preposition inside the ellipsis run of several
# "fill in the blank" lines in the prepositions worksheet.
#
# wdFindContinue = 1 ; wdReplaceOne = 1 ; wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceOne   = 1

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 6: "- beware ……………… if you do not want to be dismissed"
#            -> "- beware …of…………… if you do not want to be dismissed"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(6).Range
$p.Find.ClearFormatting()
$p.Find.Execute("beware ………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "beware …of……………", $wdReplaceOne)

# ---------------------------------------------------------------------------
# Paragraph 7: "- Getting sacked may collide ……………………"
#            -> "- Getting sacked may collide …with…………………"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(7).Range
$p.Find.ClearFormatting()
$p.Find.Execute(" ……………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, " …with…………………", $wdReplaceOne)

# ---------------------------------------------------------------------------
# Paragraph 8: "- they usually conceal reasons …………………… "
#            -> "- they usually conceal reasons …for………………… "
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(8).Range
$p.Find.ClearFormatting()
$p.Find.Execute("……………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "…for…………………", $wdReplaceOne)

# ---------------------------------------------------------------------------
# Paragraph 9: "- worker doesn't correspond ………………"
#            -> "- worker doesn't correspond …to……………"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9).Range
$p.Find.ClearFormatting()
$p.Find.Execute(" ………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, " …to……………", $wdReplaceOne)

# ---------------------------------------------------------------------------
# Paragraph 10: "- dispose …………………unnecessary workers"
#             -> "- dispose ……of……………unnecessary workers"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(10).Range
$p.Find.ClearFormatting()
$p.Find.Execute("dispose …………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "dispose ……of……………", $wdReplaceOne)

# ---------------------------------------------------------------------------
# Paragraph 13: "- object ………………………notice"
#             -> "- object ……to…………………notice"   ("to" is NOT bold, unlike the
#                other revealed prepositions above)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(13).Range
$p.Find.ClearFormatting()
$p.Find.Execute("object ………………………", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "object ……", $wdReplaceOne)

$p = $d.Paragraphs.Item(13).Range
$p.Find.ClearFormatting()
$found = $p.Find.Execute("object ……")
$insertAt = $d.Range($p.End, $p.End)
$insertAt.InsertAfter("to…………………")

$toRange = $d.Range($p.End, $p.End + 2)
$toRange.Font.Bold = 0

$dotsRange = $d.Range($p.End + 2, $p.End + 2 + 7)
$dotsRange.Font.Bold = 1
